# Add two new variables/columns to the "Data" sheet: Num_Books and Fav_Genre
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Headers (bold, like the existing Height/Weight/Sex headers) ---
$ws.Range("D1").Value = "Num_Books"
$ws.Range("E1").Value = "Fav_Genre"
$ws.Range("D1:E1").Font.Bold = $true

# --- Row data for Num_Books (D) and Fav_Genre (E) ---
$numBooks = @{
  2 = 2;  3 = 6;  4 = 1;  5 = 33; 6 = 9;
  7 = 0;  8 = 8;  9 = 20; 10 = 22; 11 = 6;
  12 = 13; 13 = 14; 14 = 11; 15 = 25
}
$genres = @{
  2 = "Romance";     3 = "Sci-Fi";      4 = "Sci-Fi";      5 = "Romance";     6 = "Dark Comedy";
  7 = "Horror";      8 = "Horror";      9 = "Sci-Fi";      10 = "Horror";     11 = "Dark Comedy";
  12 = "Sci-Fi";     13 = "Dark Comedy"; 14 = "Romance";    15 = "Sci-Fi"
}

# Write rows so that new distinct genre strings are first introduced in the
# order Sci-Fi, Horror, Romance, Dark Comedy (row 3, 7, 2, 6), then fill in
# the remaining rows in natural order.
$rowOrder = @(3, 7, 2, 6, 4, 5, 8, 9, 10, 11, 12, 13, 14, 15)
foreach ($r in $rowOrder) {
    $ws.Cells.Item($r, 4).Value = $numBooks[$r]
    $ws.Cells.Item($r, 5).Value = $genres[$r]
}

# --- View / selection tweaks from the diff ---
$ws.Range("M6").Select() | Out-Null

# --- Page setup (portrait) ---
$ws.PageSetup.Orientation = 1
